$d = $word.ActiveDocument
$vtab = [char]11

# 1. Insert a manual line break in the "Programa" paragraph, right after
#    "Função delta. " and before "Equação de Laplace". We select the
#    *entire* run's text (not just a substring) and reassign it with the
#    manual line-break character embedded, so Word emits a proper
#    xml:space="preserve" <w:t> before the <w:br/>.
$old1 = "Funções de uma variável complexa: séries infinitas, funções analíticas, condições de Cauchy Riemann, integrais de contorno, teorema de Cauchy, teorema dos resíduos, Função delta. Equação de Laplace, equação da difusão (do calor), equação de ondas (corda vibrante); Série de Fourier, Transformadas Integrais de Fourier e Laplace. Funções especiais: Polinômios de Legendre, Harmônicos Esféricos, Funções de Bessel."
$rng1 = $d.Content
$found1 = $rng1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "Could not find text for change 1" }
$splitAt = 179
$new1 = $old1.Substring(0, $splitAt) + $vtab + $old1.Substring($splitAt)
$rng1.Text = $new1

# 2. Insert manual line breaks before each bullet item in the Bibliografia
#    paragraph, turning the single run of text into four lines. Again we
#    select the entire run's text and reassign it in one shot.
$old2 = "•ARFKEN, G. and WEBER, H. J. Mathematical Methods for Physicists.•BROWN, JAMES W. and CHURCHILL, RUEL V., Complex Variables and Applications, Mc Graw Hill Higher Education, 7a. ed.• BUTKOV, Eugene. Física Matemática.•BELLANDI FILHO,J., Funções Especiais, Ed. Papirus, 1985."
$new2 = "•ARFKEN, G. and WEBER, H. J. Mathematical Methods for Physicists.•BROWN, JAMES W. and CHURCHILL, RUEL V., Complex Variables and Applications, Mc Graw Hill Higher Education, 7a. ed.• BUTKOV, Eugene. Física Matemática.•BELLANDI FILHO,J., Funções Especiais, Ed. Papirus, 1985."
$rng2 = $d.Content
$found2 = $rng2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Could not find text for change 2" }
$rng2.Text = $new2
